$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 98

# Columns A and D hold text that looks numeric/date-like (a date string and a
# zero-padded week number); prefix with a leading apostrophe so Excel stores
# them as literal text instead of auto-converting to a date serial / number.
$ws.Cells.Item($row, 1).Value = "'2025-03-03"
$ws.Cells.Item($row, 2).Value = "12:44:01"
$ws.Cells.Item($row, 3).Value = "Monday"
$ws.Cells.Item($row, 4).Value = "'09"

$ws.Cells.Item($row, 5).Value = 131912
$ws.Cells.Item($row, 6).Value = 142559
$ws.Cells.Item($row, 7).Value = 172149
$ws.Cells.Item($row, 8).Value = 159022
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 147191
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 194058
$ws.Cells.Item($row, 14).Value = 115378
$ws.Cells.Item($row, 15).Value = 46364
$ws.Cells.Item($row, 16).Value = 29762
$ws.Cells.Item($row, 17).Value = 70497
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 51002
$ws.Cells.Item($row, 20).Value = -1
